# Refactor BETA_LIGHT array for improved readability in lightmap.fsh
# (commit message is a red herring relative to the actual spreadsheet content;
#  this updates the changelog sheet with a new release entry.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new changelog entries at the bottom of the list (rows 11-12),
# and update the "current" row (row 2): Minecraft Version / Pack Version / Build.
$ws.Range("A11").Value = "Changed dried ghast textures"
$ws.Range("B2").Value = "10.0.9"
$ws.Range("A12").Value = "Released Classic Reimagined 11 (WIP) "
$ws.Range("A2").Value = "1.21.7"
$ws.Range("C2").Value = 8255

# Move the active-cell selection to E11 (matches the saved selection state)
$ws.Range("E11").Select()
